$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 2977.84
$wsSummary.Range("E2").Value = 7022.16
$wsSummary.Range("F2").Value = 46.82
$wsSummary.Range("A3").Value = 512.05999999999995
$wsSummary.Range("E3").Value = 279.52999999999997
$wsSummary.Range("A5").Value = 0.59
$wsSummary.Range("B5").Value = 0.23
$wsSummary.Range("E5").Value = 0.36

# ---------------------------------------------------------------------------
# Original Schedule sheet
# ---------------------------------------------------------------------------
$wsOriginal = $wb.Worksheets.Item("Original Schedule")
$wsOriginal.Range("F8").Value = 0
$wsOriginal.Range("G8").Value = 1055.05
$wsOriginal.Range("F9").Value = 0.36
$wsOriginal.Range("G9").Value = 1055.4100000000001

# ---------------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------------
$wsRepayment = $wb.Worksheets.Item("Repayment schedule")
$wsRepayment.Range("J8").Value = 0
$wsRepayment.Range("K8").Value = 1055.05
$wsRepayment.Range("Q8").Value = 55.1
$wsRepayment.Range("F9").Value = 983.56
$wsRepayment.Range("G9").Value = 5983.5
$wsRepayment.Range("H9").Value = 71.489999999999995
$wsRepayment.Range("J9").Value = 0.36
$wsRepayment.Range("K9").Value = 1055.4100000000001
$wsRepayment.Range("Q9").Value = 1055.4100000000001
$wsRepayment.Range("G10").Value = 4987.45
$wsRepayment.Range("G11").Value = 3983.22
$wsRepayment.Range("G12").Value = 2967.44
$wsRepayment.Range("G13").Value = 1942.62
$wsRepayment.Range("G14").Value = 907.36
$wsRepayment.Range("F15").Value = 907.36
$wsRepayment.Range("K15").Value = 916.29
$wsRepayment.Range("Q15").Value = 916.29

# ---------------------------------------------------------------------------
# Transactions sheet
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 591
$wsTransactions.Range("E2").Value = 71.47
$wsTransactions.Range("J2").Value = 7084.41
$wsTransactions.Range("A3").Value = 590
$wsTransactions.Range("A4").Value = 584
$wsTransactions.Range("E4").Value = 20.37
$wsTransactions.Range("A5").Value = 589
$wsTransactions.Range("E5").Value = 39.65
$wsTransactions.Range("J5").Value = 7012.94
$wsTransactions.Range("A6").Value = 588
$wsTransactions.Range("A7").Value = 582
$wsTransactions.Range("E7").Value = 24.86
$wsTransactions.Range("A8").Value = 586
$wsTransactions.Range("A9").Value = 587
$wsTransactions.Range("A10").Value = 576
$wsTransactions.Range("A11").Value = 575
$wsTransactions.Range("A12").Value = 574
$wsTransactions.Range("A13").Value = 567
$wsTransactions.Range("A14").Value = 560
$wsTransactions.Range("A15").Value = 559
$wsTransactions.Range("A16").Value = 558
$wsTransactions.Range("A17").Value = 549
$wsTransactions.Range("A18").Value = 548
$wsTransactions.Range("A19").Value = 547
$wsTransactions.Range("A20").Value = 536
$wsTransactions.Range("A21").Value = 535
$wsTransactions.Range("A22").Value = 534

# ---------------------------------------------------------------------------
# Selection / active-cell bookkeeping per sheet, matching the author's last
# on-screen state. Each sheet keeps its own remembered selection; the final
# Activate()+Select() determines the workbook's active sheet/tab.
# ---------------------------------------------------------------------------
$wsSummary.Activate() | Out-Null
$wsSummary.Range("C9").Select() | Out-Null

$wsOriginal.Activate() | Out-Null
$wsOriginal.Range("J12").Select() | Out-Null

$wsRepayment.Activate() | Out-Null
$wsRepayment.Range("K4").Select() | Out-Null

$wsChargesTab = $wb.Worksheets.Item("ChargesTab")
$wsChargesTab.Activate() | Out-Null
$wsChargesTab.Range("D4").Select() | Out-Null

$wsTransactions.Activate() | Out-Null
$wsTransactions.Range("G5").Select() | Out-Null
